# New trading account updates: position sizes (USD) for each holding, column E

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("E3").Value = 121360
$ws.Range("E4").Value = 177400
$ws.Range("E6").Value = 182710
$ws.Range("E8").Value = 100000
$ws.Range("E11").Value = 100000
$ws.Range("E19").Value = 55000
$ws.Range("E20").Value = 50000
$ws.Range("E25").Formula = "=329060/6"
$ws.Range("E26").Value = 50000
$ws.Range("E27").Value = 44000
$ws.Range("E29").Value = 35000
$ws.Range("E35").Value = 31000
$ws.Range("E36").Value = 31000
$ws.Range("E38").Value = 27000
$ws.Range("E39").Value = 26000

# Update the view: reset zoom to 100% and re-anchor the frozen panes/
# selection to the top of the sheet instead of the bottom.
$ws.Activate()
[void]$ws.Range("C3").Select()

$win = $excel.ActiveWindow
$win.Zoom = 100
$win.FreezePanes = $true

[void]$ws.Range("E19").Select()
